# Generate Report for Archive
#
# The localization-status report lists tracked files in ascending order by
# their (re)generation time. Two files - "e329ed08-..." and "e4b128cf-..." -
# were regenerated and the report row order now needs to swap so row 4 holds
# "e4b128cf-..." and row 5 holds "e329ed08-...". This applies to the
# "Overview" summary sheet as well as the per-locale "zh-cn" / "de-de" detail
# sheets. The swap is done in place (cell values only); row 4 and row 5 keep
# their original formatting/styles.

$wb = $excel.ActiveWorkbook

function Swap-RowValues($ws, [string[]]$cols, [int]$row1, [int]$row2) {
    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

# ---------------------------------------------------------------------
# 1) Swap row 4 / row 5 content on all three sheets.
# ---------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
Swap-RowValues $wsOverview @("A","B","C") 4 5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Swap-RowValues $wsZhCn @("A","B","C","D") 4 5

$wsDeDe = $wb.Worksheets.Item("de-de")
Swap-RowValues $wsDeDe @("A","B","C","D") 4 5

# ---------------------------------------------------------------------
# 2) Rebuild the hyperlinks on each sheet so the link text follows the row
#    it is now in (the link targets themselves stay anchored to their
#    original row/rId, matching how Excel preserves the Hyperlinks
#    collection order across a content-only edit).
# ---------------------------------------------------------------------

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c053265932b0b14eeb73e17b5af893534456ddd3/e2e/3afcb3a5-4980-43cb-9abd-59c8cdfef388.md", "", "", $wsOverview.Range("A2").Value2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c053265932b0b14eeb73e17b5af893534456ddd3/e2e/7e85abdc-023e-4001-a7e5-cfc2112e0687.md", "", "", $wsOverview.Range("A3").Value2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/20482266072b6d46ccc13484d04f18ba9d5dcdc6/e2e/e329ed08-4084-4799-9ae8-3c26ba335479.md", "", "", $wsOverview.Range("A4").Value2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/445f219e12c9d5c91a5cb64d7f2814070de222a7/e2e/e4b128cf-6e74-487c-a63b-63c0b070a7ae.md", "", "", $wsOverview.Range("A5").Value2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/20482266072b6d46ccc13484d04f18ba9d5dcdc6/.localization-config", "", "", $wsOverview.Range("A6").Value2)

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c053265932b0b14eeb73e17b5af893534456ddd3/e2e/3afcb3a5-4980-43cb-9abd-59c8cdfef388.md", "", "", $wsZhCn.Range("A2").Value2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bc65390bf6b2ee9cb3c445a7a5d1093d58dcc930/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/3afcb3a5-4980-43cb-9abd-59c8cdfef388.bb8e4420bf044dca9fa51faa7f8a0bf0de3c07ee.zh-cn.xlf", "", "", $wsZhCn.Range("C2").Value2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c053265932b0b14eeb73e17b5af893534456ddd3/e2e/7e85abdc-023e-4001-a7e5-cfc2112e0687.md", "", "", $wsZhCn.Range("A3").Value2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bc65390bf6b2ee9cb3c445a7a5d1093d58dcc930/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/7e85abdc-023e-4001-a7e5-cfc2112e0687.759644fbffe3e49ac7dffe277ff7e4f735d6b79d.zh-cn.xlf", "", "", $wsZhCn.Range("C3").Value2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/20482266072b6d46ccc13484d04f18ba9d5dcdc6/e2e/e329ed08-4084-4799-9ae8-3c26ba335479.md", "", "", $wsZhCn.Range("A4").Value2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/954253aaae46061a13120d7ea8b314eaefca7b0b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/e329ed08-4084-4799-9ae8-3c26ba335479.39805055960258112f342c2ed773a55630a474f6.zh-cn.xlf", "", "", $wsZhCn.Range("C4").Value2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/445f219e12c9d5c91a5cb64d7f2814070de222a7/e2e/e4b128cf-6e74-487c-a63b-63c0b070a7ae.md", "", "", $wsZhCn.Range("A5").Value2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/520a77e65ba1ec0642f081e505ee45216c360acc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/e4b128cf-6e74-487c-a63b-63c0b070a7ae.bccdea2384581b097eadb01ee47e5801e8184cb5.zh-cn.xlf", "", "", $wsZhCn.Range("C5").Value2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/20482266072b6d46ccc13484d04f18ba9d5dcdc6/.localization-config", "", "", $wsZhCn.Range("A6").Value2)

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c053265932b0b14eeb73e17b5af893534456ddd3/e2e/3afcb3a5-4980-43cb-9abd-59c8cdfef388.md", "", "", $wsDeDe.Range("A2").Value2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/77b40c8bca0b243664b31f667874eb161e523a3e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/3afcb3a5-4980-43cb-9abd-59c8cdfef388.bb8e4420bf044dca9fa51faa7f8a0bf0de3c07ee.de-de.xlf", "", "", $wsDeDe.Range("C2").Value2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c053265932b0b14eeb73e17b5af893534456ddd3/e2e/7e85abdc-023e-4001-a7e5-cfc2112e0687.md", "", "", $wsDeDe.Range("A3").Value2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/77b40c8bca0b243664b31f667874eb161e523a3e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/7e85abdc-023e-4001-a7e5-cfc2112e0687.759644fbffe3e49ac7dffe277ff7e4f735d6b79d.de-de.xlf", "", "", $wsDeDe.Range("C3").Value2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/20482266072b6d46ccc13484d04f18ba9d5dcdc6/e2e/e329ed08-4084-4799-9ae8-3c26ba335479.md", "", "", $wsDeDe.Range("A4").Value2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90fe93f6f8128e8bd83a3553fcaf082912293ef5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/e329ed08-4084-4799-9ae8-3c26ba335479.39805055960258112f342c2ed773a55630a474f6.de-de.xlf", "", "", $wsDeDe.Range("C4").Value2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/445f219e12c9d5c91a5cb64d7f2814070de222a7/e2e/e4b128cf-6e74-487c-a63b-63c0b070a7ae.md", "", "", $wsDeDe.Range("A5").Value2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e16b7f63d44916cbb46ce63ec7293c96718ef4af/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/e4b128cf-6e74-487c-a63b-63c0b070a7ae.bccdea2384581b097eadb01ee47e5801e8184cb5.de-de.xlf", "", "", $wsDeDe.Range("C5").Value2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/20482266072b6d46ccc13484d04f18ba9d5dcdc6/.localization-config", "", "", $wsDeDe.Range("A6").Value2)
